$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update the "data" sheet's time_taken (F) column with refreshed query timestamps ---
$newTimes = @(
    "2021-10-05 14:19:51.988075",
    "2021-10-05 14:19:51.988083",
    "2021-10-05 14:19:51.988086",
    "2021-10-05 14:19:51.988089",
    "2021-10-05 14:19:51.988092",
    "2021-10-05 14:19:51.988094",
    "2021-10-05 14:19:51.988097",
    "2021-10-05 14:19:51.988099",
    "2021-10-05 14:19:51.988102",
    "2021-10-05 14:19:51.988105",
    "2021-10-05 14:19:51.988107",
    "2021-10-05 14:19:51.988110",
    "2021-10-05 14:19:51.988112",
    "2021-10-05 14:19:51.988115",
    "2021-10-05 14:19:51.988117",
    "2021-10-05 14:19:51.988120",
    "2021-10-05 14:19:51.988122",
    "2021-10-05 14:19:51.988125",
    "2021-10-05 14:19:51.988128",
    "2021-10-05 14:19:51.988131",
    "2021-10-05 14:19:51.988133",
    "2021-10-05 14:19:51.988136",
    "2021-10-05 14:19:51.988138",
    "2021-10-05 14:19:51.988141",
    "2021-10-05 14:19:51.988144",
    "2021-10-05 14:19:51.988146",
    "2021-10-05 14:19:51.988149",
    "2021-10-05 14:19:51.988151",
    "2021-10-05 14:19:51.988154",
    "2021-10-05 14:19:51.988156",
    "2021-10-05 14:19:51.988159",
    "2021-10-05 14:19:51.988161",
    "2021-10-05 14:19:51.988164",
    "2021-10-05 14:19:51.988167",
    "2021-10-05 14:19:51.988169",
    "2021-10-05 14:19:51.988172",
    "2021-10-05 14:19:51.988174",
    "2021-10-05 14:19:51.988177",
    "2021-10-05 14:19:51.988179",
    "2021-10-05 14:19:51.988182",
    "2021-10-05 14:19:51.988185",
    "2021-10-05 14:19:51.988188",
    "2021-10-05 14:19:51.988190",
    "2021-10-05 14:19:51.988193",
    "2021-10-05 14:19:51.988195",
    "2021-10-05 14:19:51.988198",
    "2021-10-05 14:19:51.988200",
    "2021-10-05 14:19:51.988203",
    "2021-10-05 14:19:51.988205",
    "2021-10-05 14:19:51.988208",
    "2021-10-05 14:19:51.988210",
    "2021-10-05 14:19:51.988213",
    "2021-10-05 14:19:51.988215",
    "2021-10-05 14:19:51.988218",
    "2021-10-05 14:19:51.988221",
    "2021-10-05 14:19:51.988223",
    "2021-10-05 14:19:51.988226",
    "2021-10-05 14:19:51.988228",
    "2021-10-05 14:19:51.988231",
    "2021-10-05 14:19:51.988233",
    "2021-10-05 14:19:51.988236",
    "2021-10-05 14:19:51.988238",
    "2021-10-05 14:19:51.988241",
    "2021-10-05 14:19:51.988243",
    "2021-10-05 14:19:51.988247",
    "2021-10-05 14:19:51.988250",
    "2021-10-05 14:19:51.988252",
    "2021-10-05 14:19:51.988255"
)
for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" sheet (placed after "data") describing the panel query ---
$meta = $wb.Worksheets.Add($null, $ws1)
$meta.Name = "metadata"

# Reuse the bold/bordered/centered header style already used on the "data" sheet
$ws1.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Cystic kidney disease"
$meta.Range("C2").Value = 283
$meta.Range("D2").Value = "'2.26"
$meta.Range("E2").Value = "2021-04-29T22:44:29.519573Z"
$meta.Range("F2").Value = "2021-10-05 14:19:51.984238"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/283/?format=json"
